# Added code for Graded Cover
# Duplicate TC_MPM_003 into a new TC_MPM_004 sheet with Graded-Cover columns.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("TC_MPM_003")

# 1. Duplicate TC_MPM_003 -> new sheet placed right after it, then rename/activate.
$ws3.Copy([System.Type]::Missing, $ws3)
$ws4 = $wb.Worksheets.Item($ws3.Index + 1)
$ws4.Name = "TC_MPM_004"
$ws4.Activate()

# 2. TC_MPM_003's own selection becomes "whole data rows selected" after the edit.
$ws3.Rows("1:2").Select()

# 3. AD2 on the new sheet switches from "Multiple Of Salary" to "Graded Cover".
$ws4.Range("AD2").Value = "Graded Cover"

# 4. Replace the old single-grade summary columns (AE:AH) with the new
#    Graded Cover block (AE:AS) -- 3 grades x (label + 4 value columns).
$ws4.Range("AE1:AH2").ClearContents()

function Set-GradedHeader {
    param($cell, $label, $gradeSuffix)
    $full = "$label`n($gradeSuffix)"
    $cell.Value = $full
    $labelLen = ($label + "`n").Length
    $tailLen = "($gradeSuffix)".Length
    $chars = $cell.Characters($labelLen + 1, $tailLen)
    $chars.Font.Color = 65535
    $chars.Font.Name = "Calibri"
}

$grades = "Grade 1", "Grade 2", "Grade 3"
$headerCols = "AF", "AG", "AH", "AI", "AK", "AL", "AM", "AN", "AP", "AQ", "AR", "AS"
$labels = "Sum Assured", "Free Cover Limit (Sum Assured)", "Minimum Cap", "Maximum Cap "
$gradeHeaderCols = "AE", "AJ", "AO"

for ($g = 0; $g -lt 3; $g++) {
    $gradeLabel = $grades[$g]
    $ws4.Range($gradeHeaderCols[$g] + "1").Value = $gradeLabel
    for ($i = 0; $i -lt 4; $i++) {
        $col = $headerCols[$g * 4 + $i]
        Set-GradedHeader -cell $ws4.Range($col + "1") -label $labels[$i] -gradeSuffix $gradeLabel
    }
}

# 5. Data row (row 2): grade name + Sum Assured / Free Cover Limit / Min Cap / Max Cap.
$dataCols = "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM", "AN", "AO", "AP", "AQ", "AR", "AS"
$dataVals = "Grade 1", "1000000", "5000000", "1000", "10000000", `
            "Grade 2", "2000000", "5000000", "1000", "10000000", `
            "Grade 3", "3000000", "5000000", "1000", "10000000"
$isNumericLooking = $false, $true, $true, $true, $true, `
                     $false, $true, $true, $true, $true, `
                     $false, $true, $true, $true, $true

for ($i = 0; $i -lt $dataCols.Length; $i++) {
    $rng = $ws4.Range($dataCols[$i] + "2")
    if ($isNumericLooking[$i]) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $dataVals[$i]
}

# 6. Row heights: header row grows slightly to fit the new wrapped headers;
#    the data row goes back to the sheet default (no explicit override).
$ws4.Rows.Item(1).RowHeight = 129
$ws4.Rows.Item(2).AutoFit()

# 7. New columns get sensible widths (closest representable values).
$colWidths = @{
    "A" = 19.666666666666668; "B" = 19.0; "C" = 19.833333333333332; "D" = 16.666666666666668;
    "E" = 12.166666666666666; "F" = 14.666666666666666; "G" = 15.5; "H" = 13.5;
    "I" = 26.666666666666668; "J" = 15.666666666666666; "K" = 10.5;
    "M" = 20.333333333333332; "N" = 13.666666666666666; "O" = 23.166666666666668; "P" = 16.666666666666668;
    "S" = 18.0; "T" = 18.833333333333332; "U" = 17.166666666666668; "V" = 18.833333333333332; "W" = 30.5;
    "Z" = 23.166666666666668; "AA" = 19.666666666666668; "AB" = 16.166666666666668;
    "AC" = 51.833333333333336; "AD" = 46.0;
    "AE" = 11.666666666666666; "AF" = 19.5; "AG" = 16.666666666666668; "AH" = 13.666666666666666;
    "AI" = 16.5; "AK" = 18.0; "AL" = 25.5; "AM" = 19.5; "AN" = 15.833333333333334;
    "AP" = 19.833333333333332; "AQ" = 25.166666666666668; "AR" = 15.166666666666666; "AS" = 18.166666666666668
}
foreach ($colLetter in $colWidths.Keys) {
    $ws4.Columns.Item($colLetter).ColumnWidth = $colWidths[$colLetter]
}

# 8. View state for the new sheet: scrolled to the Graded Cover block, cursor on AM6.
$excel.ActiveWindow.ScrollColumn = 31
$ws4.Range("AM6").Select()

# 9. Page setup (portrait) for the new sheet, matching its siblings.
$ws4.PageSetup.Orientation = 1
